$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.792.48"
$ws.Range("E2").Value = "  -1.81%  "
$ws.Range("D3").Value = "1.875.81"
$ws.Range("E3").Value = "  -2.50%  "
$ws.Range("E4").Value = "  -0.73%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.85"
$ws.Range("E5").Value = "  -2.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.684"
$ws.Range("E6").Value = "  -7.55%  "
$ws.Range("E7").Value = "  -0.76%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.87"
$ws.Range("E8").Value = "  +2.75%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.346"
$ws.Range("E9").Value = "  -3.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "51.12"
$ws.Range("E10").Value = "  -2.90%  "
$ws.Range("E11").Value = "  -0.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0971"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "12.92"
$ws.Range("E13").Value = "  +1.64%  "
$ws.Range("D14").Value = "2.149.41"
$ws.Range("E14").Value = "  -2.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.716"
$ws.Range("E15").Value = "  -0.20%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.89"
$ws.Range("E16").Value = "  -0.32%  "
$ws.Range("D17").Value = "1.882.66"
$ws.Range("E17").Value = "  -2.49%  "
$ws.Range("D18").Value = "34.823.37"
$ws.Range("E18").Value = "  -1.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.82"
$ws.Range("E19").Value = "  -0.96%  "
$ws.Range("D20").Value = "0.0₃0820"
$ws.Range("E20").Value = "  -2.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "243.45"
$ws.Range("E21").Value = "  +0.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.73"
$ws.Range("E22").Value = "  -2.54%  "
$ws.Range("E23").Value = "  -2.94%  "
$ws.Range("E24").Value = "  -0.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.45"
$ws.Range("E25").Value = "  +4.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.24"
$ws.Range("E26").Value = "  -4.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.18"
$ws.Range("E27").Value = "  -1.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.38"
$ws.Range("E28").Value = "  -3.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.20"
$ws.Range("E29").Value = "  -3.58%  "
$ws.Range("E30").Value = "  -6.71%  "
$ws.Range("D31").Value = "4.128.43"
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("E32").Value = "  +1.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.23"
$ws.Range("E33").Value = "  -2.96%  "
$ws.Range("E34").Value = "  -0.62%  "
$ws.Range("E35").Value = "  -2.34%  "
$ws.Range("E36").Value = "  -0.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.828"
$ws.Range("E37").Value = "  -9.40%  "
$ws.Range("E38").Value = "  -2.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.52"
$ws.Range("E39").Value = "  -23.36%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "97.71"
$ws.Range("E40").Value = "  -1.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.03"
$ws.Range("E41").Value = "  -1.95%  "
$ws.Range("E42").Value = "  +0.78%  "
$ws.Range("E43").Value = "  +0.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.09"
$ws.Range("E44").Value = "  -4.47%  "
$ws.Range("D45").Value = "1.284.79"
$ws.Range("E45").Value = "  -4.78%  "
$ws.Range("E46").Value = "  -6.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0804"
$ws.Range("E47").Value = "  +9.88%  "
$ws.Range("E48").Value = "  -1.10%  "
$ws.Range("E49").Value = "  -2.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "12.10"
$ws.Range("E50").Value = "  +5.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.43"
$ws.Range("E51").Value = "  -4.05%  "
